$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K financial data to E:L
# (adds a new most-recent fiscal year/period as the first data column)
$ws.Range("D1").EntireColumn.Insert()

# Copy number formats/styles from column E (the old column D, now shifted) into new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth()

# Populate the new column D with the newest period's figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1679700
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 582000
$ws.Range("D17").Value = 696800
$ws.Range("D18").Value = 982900
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 1564900
$ws.Range("D22").Value = 342700
$ws.Range("D23").Value = 640100
$ws.Range("D24").Value = 129300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 510800
$ws.Range("D27").Value = 510800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 510800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 510800
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 300100
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 15707100
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 22900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 18481800
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 1372700
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 11538900
$ws.Range("D62").Value = 643800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 13674900
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 2331600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 4806900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 510800
$ws.Range("D83").Value = 582000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1254100
$ws.Range("D91").Value = -2800100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -3384800
$ws.Range("D96").Value = -41600
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 2145400
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 14700

# Minor restatements to the two most recent periods in a couple of rows
# (rows 100 and 102 - shifted columns E and F get slightly revised figures)
$ws.Range("E100").Value = 1101700
$ws.Range("F100").Value = 1103000
$ws.Range("E102").Value = 17500
$ws.Range("F102").Value = 117600
